$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All cells in columns B-E of this sheet hold text data (coin name, link,
# price, volume change). Column D (Price) values sometimes look like plain
# decimal numbers (e.g. "214.70"); format those cells as Text first so Excel
# keeps them as strings instead of silently converting them to numeric values.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.870.86'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.667.82'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.70'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.515'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.63'
$ws.Range('E8').Value = '  +3.73%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0878'
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.903.85'
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.659.34'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.93'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '252.32'
$ws.Range('E17').Value = '  +7.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.854.98'
$ws.Range('E18').Value = '  +2.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0732'
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('E20').Value = '  -4.21%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.47'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.39'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.81'
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.24'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('E30').Value = '  +5.68%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.14'
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.416.99'
$ws.Range('E34').Value = '  -8.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('E35').Value = '  -5.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.38'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  -1.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.579'
$ws.Range('E38').Value = '  -4.62%  '
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('E40').Value = '  -2.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.45'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.811.31'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.38'
$ws.Range('E45').Value = '  -7.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.790'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  +4.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.74'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('E50').Value = '  -2.14%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0511'
$ws.Range('E51').Value = '  +0.08%  '
